$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames (row 1) ---
# LEGENDA -> legenda
# area_km2 -> area
# nome -> nome (unchanged)
# area_km2_1 -> area_km2
$ws.Range("A1").Value = "legenda"
$ws.Range("B1").Value = "area"
$ws.Range("C1").Value = "nome"
$ws.Range("D1").Value = "area_km2"

# --- Updated area_km2 values (column B) ---
$ws.Range("B2").Value  = 15218.4500735
$ws.Range("B3").Value  = 2071.24351554
$ws.Range("B4").Value  = 11536.911461
$ws.Range("B5").Value  = 1675.69051521
$ws.Range("B6").Value  = 2352.8960852
$ws.Range("B7").Value  = 4150.10357051
$ws.Range("B8").Value  = 9619.196072950001
$ws.Range("B9").Value  = 7226.80132027
$ws.Range("B10").Value = 1893.2236196
$ws.Range("B11").Value = 947.873678947
$ws.Range("B12").Value = 12098.4580798
$ws.Range("B13").Value = 2572.3296401
$ws.Range("B14").Value = 338.271731251
$ws.Range("B15").Value = 2141.05750381
$ws.Range("B16").Value = 19994.0316874
$ws.Range("B17").Value = 4726.87460522
$ws.Range("B18").Value = 10973.9920241
$ws.Range("B19").Value = 5534.13484415
$ws.Range("B20").Value = 9920.822482580001
$ws.Range("B21").Value = 10093.9755614
$ws.Range("B22").Value = 4161.26554807
$ws.Range("B23").Value = 9863.62618297
$ws.Range("B24").Value = 8445.11739858
$ws.Range("B25").Value = 1373.4916239
$ws.Range("B26").Value = 15457.276101
